$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get purely-numeric-looking replacement text (e.g. "298.80", "0.492");
# force Text format first so Excel keeps them as strings instead of auto-converting to numbers,
# matching the original inlineStr/text storage of column D.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = '42.177.43'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '2.266.60'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '298.80'
$ws.Range("E5").Value = '  -1.94%  '
$ws.Range("D6").Value = '95.53'
$ws.Range("E6").Value = '  -5.00%  '
$ws.Range("E7").Value = '  -2.43%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").Value = '33.32'
$ws.Range("E10").Value = '  -3.46%  '
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '47.98'
$ws.Range("E12").Value = '  -8.04%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").Value = '6.65'
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").Value = '2.619.86'
$ws.Range("E15").Value = '  -2.49%  '
$ws.Range("D16").Value = '15.51'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '2.276.32'
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").Value = '42.098.97'
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("D20").Value = '11.66'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("E22").Value = '  -2.83%  '
$ws.Range("D23").Value = '66.53'
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("D24").Value = '234.25'
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("E25").Value = '  -1.84%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -3.18%  '
$ws.Range("D28").Value = '23.95'
$ws.Range("E28").Value = '  -5.71%  '
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("D30").Value = '168.40'
$ws.Range("E30").Value = '  +5.09%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '9.17'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '33.68'
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '4.89'
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").Value = '4.45'
$ws.Range("E35").Value = '  -3.31%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.33'
$ws.Range("E36").Value = '  -4.84%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").Value = '16.47'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("D38").Value = '0.0685'
$ws.Range("E38").Value = '  -4.90%  '
$ws.Range("E39").Value = '  -3.86%  '
$ws.Range("D40").Value = '0.0984'
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("E42").Value = '  -5.73%  '
$ws.Range("E43").Value = '  -4.71%  '
$ws.Range("D44").Value = '1.957.05'
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("E45").Value = '  -1.88%  '
$ws.Range("D46").Value = '17.39'
$ws.Range("E46").Value = '  -7.26%  '
$ws.Range("D47").Value = '9.53'
$ws.Range("E47").Value = '  -6.21%  '
$ws.Range("E48").Value = '  -3.97%  '
$ws.Range("D49").Value = '2.492.10'
$ws.Range("E50").Value = '  -5.81%  '
$ws.Range("E51").Value = '  -3.17%  '
